$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update 想去人数 (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 6625
$ws1.Range("F6").Value = 5402
$ws1.Range("F12").Value = 35

# Sheet "全部类型" (All Types) - same events, update matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 6625
$ws4.Range("F6").Value = 5402
$ws4.Range("F14").Value = 35
